# fix: avoid saving first excel row
#
# The workbook's first row was a merged "Table 1" title banner sitting above
# the real header row. That banner row should not have been saved/exported;
# remove it so the real header row becomes row 1 and the frozen-pane view
# (which was anchored around the old title/header rows) no longer applies.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the stray first row (the "Table 1" banner) - remaining rows shift up.
$ws.Rows("1:1").Delete()

# The old frozen pane was anchored to the banner/header rows (topLeftCell
# B3, 1 column x 2 rows frozen). With the banner row gone that split no
# longer makes sense, so drop the freeze entirely.
$excel.ActiveWindow.FreezePanes = $false
